$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.092.06"

$ws.Range("D3").Value = "3.866.13"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "698.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.11%  "

$ws.Range("D7").Value = "3.863.53"
$ws.Range("E7").Value = "  +1.34%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  +1.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.39%  "

$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("E13").Value = "  +5.00%  "

$ws.Range("E14").Value = "  +1.17%  "

$ws.Range("D15").Value = "4.518.19"
$ws.Range("E15").Value = "  +1.42%  "

$ws.Range("D16").Value = "3.861.15"
$ws.Range("E16").Value = "  +1.32%  "

$ws.Range("D17").Value = "71.164.11"
$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.77%  "

$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "495.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.89%  "

$ws.Range("E23").Value = "  +1.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.26%  "

$ws.Range("E25").Value = "  +1.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.83%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("E28").Value = "  +1.47%  "

$ws.Range("D29").Value = "4.007.02"
$ws.Range("E29").Value = "  +1.11%  "

$ws.Range("E30").Value = "  +7.97%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.07%  "

$ws.Range("E33").Value = "  -0.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.51%  "

$ws.Range("E35").Value = "  -1.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.74%  "

$ws.Range("D37").Value = "3.817.69"
$ws.Range("E37").Value = "  +1.46%  "

$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.104"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.87%  "

$ws.Range("E40").Value = "  +11.79%  "

$ws.Range("E41").Value = "  +1.97%  "

$ws.Range("E42").Value = "  -0.15%  "

$ws.Range("E43").Value = "  +7.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.80%  "

$ws.Range("E47").Value = "  +4.33%  "

$ws.Range("E48").Value = "  +0.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.43%  "

$ws.Range("E50").Value = "  +1.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "417.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.98%  "
